$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# AppControl sheet: the "Email ID" distribution list (B25) grows from
# two addresses to the full team list, and it stops being a mailto:
# hyperlink (plain text instead).
# ------------------------------------------------------------------
$appControl = $wb.Worksheets.Item("AppControl")

# Drop the mailto: hyperlink that lived on B25.
$appControl.Hyperlinks.Delete()

# The cell's style still points at the built-in "Hyperlink" cell style;
# remove that named style now that nothing should look like a link.
$wb.Styles.Item("Hyperlink").Delete()

$b25 = $appControl.Range("B25")
$b25.Value = "nilesh@zestiot.io, amit@zestiot.io, sushanto@zestiot.io, sudhir@zestiot.io, Krishna@zestiot.io, anantwar@zestiot.io, shrikant@zestiot.io, aman@zestiot.io, rohan@zestiot.io, leadership@enhops.com, chiranjeevi@zestiot.io, stiyyagura@enhops.com, pdwadasi@enhops.com, rbuddha@enhops.com, rchiluka@enhops.com, smunnangi@enhops.com, nishanth@zestiot.io, hmanthena@enhops.com, mpyla@enhops.com"

# Plain (non-hyperlink) look: regular black Segoe UI text, no underline.
$b25.Font.Name = "Segoe UI"
$b25.Font.Underline = $false
$b25.Font.ColorIndex = 1

# The much longer text now needs a taller row to keep wrapping nicely.
$appControl.Rows.Item(25).RowHeight = 247.5

# Leave the selection where the edit happened.
$appControl.Range("B25").Select()

# ------------------------------------------------------------------
# smoke sheet: the whole suite is switched on - Run Flag N -> Y for
# every test row (B3:B18).
# ------------------------------------------------------------------
$smoke = $wb.Worksheets.Item("smoke")
$smoke.Range("B3:B18").Value = "Y"

$smoke.Range("B19").Select()
